$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.733.08'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '2.712.53'
$ws.Range('E3').Value = '  +2.41%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.89'
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '163.18'
$ws.Range('E6').Value = '  +3.79%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = '2.709.88'
$ws.Range('E9').Value = '  +2.36%  '
$ws.Range('E10').Value = '  -0.39%  '
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.361'
$ws.Range('E13').Value = '  +2.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.50'
$ws.Range('E14').Value = '  +1.53%  '
$ws.Range('D15').Value = '3.194.11'
$ws.Range('E15').Value = '  +2.03%  '
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').Value = '68.717.30'
$ws.Range('E17').Value = '  +0.65%  '
$ws.Range('D18').Value = '2.719.14'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.84'
$ws.Range('E19').Value = '  +4.02%  '
$ws.Range('E20').Value = '  +4.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '365.23'
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.54'
$ws.Range('E22').Value = '  +2.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.93'
$ws.Range('E23').Value = '  +2.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.12'
$ws.Range('E24').Value = '  +2.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.04'
$ws.Range('E25').Value = '  -1.39%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  +1.65%  '
$ws.Range('D28').Value = '2.840.70'
$ws.Range('E28').Value = '  +2.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000105'
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '594.99'
$ws.Range('E30').Value = '  +6.29%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.31'
$ws.Range('E32').Value = '  +3.15%  '
$ws.Range('E33').Value = '  +3.27%  '
$ws.Range('E34').Value = '  +4.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.132'
$ws.Range('E35').Value = '  +3.08%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.87'
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '160.20'
$ws.Range('E39').Value = '  -0.43%  '
$ws.Range('E40').Value = '  +2.46%  '
$ws.Range('E41').Value = '  +2.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.45'
$ws.Range('E42').Value = '  +2.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.70'
$ws.Range('E43').Value = '  +3.24%  '
$ws.Range('E44').Value = '  +1.21%  '
$ws.Range('D45').Value = '0.0₆0318'
$ws.Range('E45').Value = '  -4.98%  '
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '158.46'
$ws.Range('E47').Value = '  -0.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.94'
$ws.Range('E48').Value = '  +5.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.79'
$ws.Range('E49').Value = '  +6.02%  '
$ws.Range('E50').Value = '  +7.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.14'
$ws.Range('E51').Value = '  +0.13%  '
